# Workbook "Hortaliza, Agricola del Norte S.A. de Arica - Cilantro"
# A new daily price observation is inserted at row 14, pushing all the
# existing observations (rows 14-106) down by one row. The row that used
# to be the last data row (106) ends up surviving as row 107 automatically
# because Excel's row-insert shifts the whole sheet down - we only need to
# populate the brand-new row 14 with its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 14:106 down to 15:107, inserting a new blank row 14.
$ws.Rows("14").Insert()

# Populate the new row 14 with the new observation.
$ws.Range("A14").Value = 1
$ws.Range("B14").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C14").Value = "Arica y Parinacota"
$ws.Range("D14").Value = 44970
$ws.Range("E14").Value = 15
$ws.Range("F14").Value = 100112040
$ws.Range("G14").Value = "Cilantro"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 300
$ws.Range("K14").Value = 2300
$ws.Range("L14").Value = 2500
$ws.Range("M14").Value = 2400
$ws.Range("N14").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("O14").Value = "Región de Arica y Parinacota"
$ws.Range("P14").Value = 1200
$ws.Range("Q14").Value = 2
$ws.Range("R14").Value = "Hortaliza"

# Make sure the date cell keeps the same date/time number format as the
# other date cells in column D.
$ws.Range("D14").NumberFormat = $ws.Range("D15").NumberFormat
